$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New H6: link to the newly added class 5 slides (set first so the shared
# string table reuses the slot being vacated by H5's old value).
$ws.Range("H6").Value = "[Slides](slides/class_5/class_5#1) [.Rmd](slides/class_5/class_5.Rmd)"

# Update H5: was the "class 4" slides link (corrected so the folder is
# slides/class_4 instead of slides/class_2).
$ws.Range("H5").Value = "[Slides](slides/class_4/class_2#1) [.Rmd](slides/class_4/class_4.Rmd) [.R](slides/class_4/class_4_taller.R)"

# Update the view: scroll so column D is the left-most visible column,
# and move the active selection to H6.
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("H6").Select()
